$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 0.7142857142857143
$ws.Range("F3").Value = 0.9767441860465116

# Row 33
$ws.Range("D33").Value = 0.7142857142857143
$ws.Range("F33").Value = 0.9767441860465116

# Row 77
$ws.Range("D77").Value = 0.8
$ws.Range("F77").Value = 0.9651162790697675

# Row 97
$ws.Range("D97").Value = 1
$ws.Range("F97").Value = 0.9767441860465116

# Row 99
$ws.Range("D99").Value = 0.4
$ws.Range("F99").Value = 0.9186046511627907
